$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.942.50'
$ws.Range('E2').Value = '  -6.12%  '
$ws.Range('D3').Value = '2.988.44'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '571.30'
$ws.Range('E5').Value = '  -3.96%  '
$ws.Range('D6').Value = '125.09'
$ws.Range('E6').Value = '  -8.95%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '2.985.69'
$ws.Range('E8').Value = '  -6.44%  '
$ws.Range('E9').Value = '  -2.58%  '
$ws.Range('E10').Value = '  -9.50%  '
$ws.Range('E11').Value = '  -5.93%  '
$ws.Range('E12').Value = '  -4.11%  '
$ws.Range('E13').Value = '  -9.77%  '
$ws.Range('E14').Value = '  -7.26%  '
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('D16').Value = '3.479.66'
$ws.Range('E16').Value = '  -6.48%  '
$ws.Range('D17').Value = '2.977.37'
$ws.Range('E17').Value = '  -6.72%  '
$ws.Range('D18').Value = '59.951.32'
$ws.Range('E18').Value = '  -6.10%  '
$ws.Range('D19').Value = '6.45'
$ws.Range('E19').Value = '  -2.28%  '
$ws.Range('D20').Value = '426.20'
$ws.Range('E20').Value = '  -8.21%  '
$ws.Range('E21').Value = '  -6.62%  '
$ws.Range('D22').Value = '0.668'
$ws.Range('E22').Value = '  -5.04%  '
$ws.Range('E23').Value = '  -8.67%  '
$ws.Range('D24').Value = '12.90'
$ws.Range('E24').Value = '  -2.82%  '
$ws.Range('E25').Value = '  -5.02%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D28').Value = '2.53'
$ws.Range('E28').Value = '  -6.29%  '
$ws.Range('E29').Value = '  -7.62%  '
$ws.Range('D30').Value = '7.25'
$ws.Range('E30').Value = '  -7.73%  '
$ws.Range('D31').Value = '6.14'
$ws.Range('E31').Value = '  -11.07%  '
$ws.Range('D32').Value = '25.15'
$ws.Range('E32').Value = '  -9.05%  '
$ws.Range('D33').Value = '0.0936'
$ws.Range('E33').Value = '  -8.46%  '
$ws.Range('E34').Value = '  -5.26%  '
$ws.Range('D35').Value = '0.928'
$ws.Range('E35').Value = '  -9.50%  '
$ws.Range('D36').Value = '50.06'
$ws.Range('E36').Value = '  -3.24%  '
$ws.Range('E37').Value = '  -17.32%  '
$ws.Range('D38').Value = '0.0₃0660'
$ws.Range('E38').Value = '  -11.28%  '
$ws.Range('D39').Value = '8.36'
$ws.Range('E39').Value = '  +1.91%  '
$ws.Range('D40').Value = '0.0353'
$ws.Range('E40').Value = '  -10.80%  '
$ws.Range('E41').Value = '  -5.76%  '
$ws.Range('D42').Value = '374.94'
$ws.Range('E42').Value = '  -5.54%  '
$ws.Range('D43').Value = '2.666.86'
$ws.Range('E43').Value = '  -5.09%  '
$ws.Range('D44').Value = '2.45'
$ws.Range('E44').Value = '  -8.86%  '
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('E46').Value = '  -8.13%  '
$ws.Range('D47').Value = '119.76'
$ws.Range('E47').Value = '  -7.31%  '
$ws.Range('E48').Value = '  -7.69%  '
$ws.Range('E49').Value = '  -4.32%  '
$ws.Range('D50').Value = '23.49'
$ws.Range('E50').Value = '  -8.67%  '
$ws.Range('E51').Value = '  -8.35%  '
